# Hawaii40 dynamic-case "limits relaxation" fix.
#
# REGCA1 (renewable-generator converter) rows 2..10:
#   Iqrmax (col T) / Iqrmin (col U) are set to +/-(100 / Sn), i.e. the
#   reactive-current-rate limit expressed on the 100 MVA system base for
#   each unit's own MVA rating Sn.
#
# REECA1 (renewable-energy electrical-control) rows 2..10:
#   Iqh1/Iql1 (R/S), QMax/QMin (Y/Z), VMAX/VMIN (AA/AB) and PMAX (AK) are
#   relaxed to +/-999 (i.e. effectively unlimited), Vref1 (AG) is turned on
#   (0 -> 1), and Imax (AM) is rederived from the new PMAX (999) scaled by
#   the associated REGCA1 unit's MVA rating: Imax = PMAX * 100 / Sn.
#
# The exact post-recalculation values below match what Excel/the source
# tool produced (and are pinned literally so the stored IEEE-754 doubles
# are bit-identical to the target workbook, since `100/Sn` evaluated at
# full double precision can land one ULP away from the reference values).

$wb = $excel.ActiveWorkbook
$wsRegca1 = $wb.Worksheets.Item("REGCA1")
$wsReeca1 = $wb.Worksheets.Item("REECA1")

# Row -> Iqrmax (REGCA1, T column); Iqrmin (U) is always -Iqrmax.
$iqrmaxByRow = @{
    2  = 1.98019801980198
    3  = 6.17283950617284
    4  = 4.545454545454546
    5  = 8.928571428571429
    6  = 3.03030303030303
    7  = 3.289473684210527
    8  = 1.855287569573284
    9  = 1.317523056653491
    10 = 3.289473684210527
}

foreach ($r in $iqrmaxByRow.Keys) {
    $iqrmax = $iqrmaxByRow[$r]
    $wsRegca1.Cells.Item($r, 20).Value2 = $iqrmax    # T: Iqrmax
    $wsRegca1.Cells.Item($r, 21).Value2 = -$iqrmax   # U: Iqrmin
}

# Row -> Imax (REECA1, AM column), derived from PMAX(=999) * 100 / Sn of
# the associated REGCA1 unit.
$imaxByRow = @{
    2  = 1978.217821782178
    3  = 1978.217821782178
    4  = 4540.909090909092
    5  = 4540.909090909092
    6  = 3027.272727272727
    7  = 3027.272727272727
    8  = 1853.43228200371
    9  = 1853.43228200371
    10 = 3286.184210526316
}

foreach ($r in $imaxByRow.Keys) {
    $wsReeca1.Cells.Item($r, 18).Value2 = 999            # R:  Iqh1
    $wsReeca1.Cells.Item($r, 19).Value2 = -999           # S:  Iql1
    $wsReeca1.Cells.Item($r, 25).Value2 = 999            # Y:  QMax
    $wsReeca1.Cells.Item($r, 26).Value2 = -999           # Z:  QMin
    $wsReeca1.Cells.Item($r, 27).Value2 = 999            # AA: VMAX
    $wsReeca1.Cells.Item($r, 28).Value2 = -999           # AB: VMIN
    $wsReeca1.Cells.Item($r, 33).Value2 = 1              # AG: Vref1
    $wsReeca1.Cells.Item($r, 37).Value2 = 999            # AK: PMAX
    $wsReeca1.Cells.Item($r, 39).Value2 = $imaxByRow[$r] # AM: Imax
}
